$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 — new car enquiry record.
# Leading "'" forces text storage for values that otherwise look numeric
# (matches how the source data is stored as text in the other rows).
$ws.Range("A3").Value = "Tata"
$ws.Range("B3").Value = "Model 2024"
$ws.Range("C3").Value = "'2020"
$ws.Range("D3").Value = "Diesel"
$ws.Range("E3").Value = "Chennai"
$ws.Range("F3").Value = "TamilNadu"
$ws.Range("G3").Value = "Chandru"
$ws.Range("H3").Value = "Chandru"
$ws.Range("I3").Value = "'7092312288"
$ws.Range("J3").Value = "chandru@gmail.com"
$ws.Range("K3").Value = "Email"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = $true
